$wb = $excel.ActiveWorkbook

# Sheet: hpi
$ws = $wb.Worksheets.Item('hpi')
$ws.Range('A2').Value = 'A long-standing history of progressive dysphagia accompanied by regurgitation of undigested food.'
$ws.Range('B2').Value = 'This finding is highly supportive as progressive dysphagia is common in esophageal dysmotility disorders, including Type 2 Achalasia, and is frequently seen in patients with CREST syndrome who develop esophageal involvement.'
$ws.Range('C2').Value = 'Lack of any history of digital color changes or vasospastic events.'
$ws.Range('D2').Value = 'The absence of Raynaud''s phenomenon, a highly specific and sensitive feature of CREST syndrome, weighs against the diagnosis despite the presence of achalasia.'
$ws.Range('A3').Value = 'Reports of achy chest discomfort and cough that worsen with meals.'
$ws.Range('B3').Value = 'These symptoms are consistent with esophageal motor dysfunction seen in Type 2 Achalasia, strengthening the association with CREST syndrome where esophageal dysmotility is a hallmark.'
$ws.Range('C3').Value = 'No mention of skin tightening or thickening on the fingers over time.'
$ws.Range('D3').Value = 'The absence of sclerodactyly, which is a key component of CREST syndrome, reduces the likelihood of the diagnosis significantly.'
$ws.Range('A4').Value = 'Patient describes phenomena consistent with Raynaud''s episodes, with color changes in fingers on cold exposure.'
$ws.Range('B4').Value = 'Raynaud''s phenomenon is one of the components of CREST syndrome and has a high specificity when present, lending strong support for the diagnosis.'
$ws.Range('C4').Value = 'Patient denies any history of calcinosis, such as palpable subcutaneous calcific deposits.'
$ws.Range('D4').Value = 'The absence of calcinosis, though not as sensitive, still detracts from a full spectrum presentation of CREST syndrome, thereby arguing against the diagnosis.'
$ws.Range('A5').Value = 'History of thickened, tight skin on the fingers (suggesting sclerodactyly).'
$ws.Range('B5').Value = 'Sclerodactyly is a classic component of CREST syndrome, and when present in the history, it significantly increases the likelihood of a CREST diagnosis with associated esophageal dysmotility.'
$ws.Range('C5').Value = 'No reports of telangiectasias on the face or hands.'
$ws.Range('D5').Value = 'Lack of telangiectasia, another highly specific qualitative marker of CREST syndrome, weakens the supporting evidence for CREST-related Type 2 Achalasia.'
$ws.Range('A6').Value = 'Episodes of small, visible dilated blood vessels on the face and hands (telangiectasias).'
$ws.Range('B6').Value = 'Telangiectasia is a highly specific sign in CREST syndrome and its presence in the history further supports the diagnosis alongside esophageal motility issues characteristic of Type 2 Achalasia.'
$ws.Range('C6').Value = 'The history is more consistent with primary gastroesophageal reflux disease (GERD) rather than a progressive motility disorder.'
$ws.Range('D6').Value = 'A history suggestive of GERD without progressive dysphagia or regurgitation typical of achalasia makes the scenario less supportive of CREST syndrome with esophageal dysmotility.'

# Sheet: hist
$ws = $wb.Worksheets.Item('hist')
$ws.Range('A2').Value = 'Past Medical History documenting Raynaud''s phenomenon'
$ws.Range('B2').Value = 'Raynaud''s phenomenon is one of the hallmark features of CREST syndrome and its frequent occurrence in these patients makes it a highly specific supporting evidence even when considered within a hist category.'
$ws.Range('C2').Value = 'Past Medical History with isolated gastroesophageal reflux disease (GERD) without any mention of connective tissue disease signs'
$ws.Range('D2').Value = 'GERD alone is common in the general population and does not specifically support CREST syndrome, reducing the likelihood ratio when no other systemic features are present.'
$ws.Range('A3').Value = 'History of sclerodactyly'
$ws.Range('B3').Value = 'Sclerodactyly, or skin thickening of the fingers, is a core component of CREST syndrome. Its presence strongly supports the diagnosis as it is rarely seen outside connective tissue disorders.'
$ws.Range('C3').Value = 'Absence of documented Raynaud''s phenomenon in the patient’s history'
$ws.Range('D3').Value = 'Without Raynaud''s phenomenon, one of the most sensitive and specific indicators of CREST syndrome is missing, strongly arguing against the full diagnosis.'
$ws.Range('A4').Value = 'Past Medical History noting calcinosis cutis'
$ws.Range('B4').Value = 'Calcinosis (calcium deposits in the skin) is a specific finding within CREST syndrome that, when present, increases the likelihood ratio favoring the diagnosis.'
$ws.Range('C4').Value = 'No history of sclerodactyly or other skin changes typical of systemic sclerosis'
$ws.Range('D4').Value = 'The absence of sclerodactyly, a key diagnostic marker of CREST syndrome, diminishes the specificity required to support the condition.'
$ws.Range('A5').Value = 'Documentation of telangiectasias'
$ws.Range('B5').Value = 'Telangiectasias are a classic feature of CREST syndrome. Their documentation provides strong, specific evidence in favor of the disease when observed in the patient''s history.'
$ws.Range('C5').Value = 'Lack of previous surgical or therapeutic interventions for esophageal dysmotility'
$ws.Range('D5').Value = 'If there is no history of treatment addressing esophageal dysmotility, then the presence of Type 2 Achalasia becomes less likely, arguing against its role in a CREST syndrome diagnosis.'
$ws.Range('A6').Value = 'Past Medical History of esophageal dysmotility management, such as prior pneumatic dilation or pharmacologic therapy for Achalasia'
$ws.Range('B6').Value = 'The treatment history for esophageal dysmotility correlates with Type 2 Achalasia, a key component of the diagnosis. This suggests symptom management reflective of the disease process.'
$ws.Range('C6').Value = 'Past Medical History free from other connective tissue disease indicators (e.g., no calcinosis or telangiectasias)'
$ws.Range('D6').Value = 'The complete absence of additional CREST-related features further lowers the pre-test probability for the syndrome, serving as strong evidence against the diagnosis.'

# Sheet: soc
$ws = $wb.Worksheets.Item('soc')
$ws.Range('A2').Value = 'Positive family history of autoimmune disorders'
$ws.Range('B2').Value = 'A family history of autoimmune diseases increases the pre-test probability for CREST syndrome, making an autoimmune connective tissue process (including CREST with associated esophageal dysmotility) more likely.'
$ws.Range('C2').Value = 'Absence of a family history of autoimmune diseases or connective tissue disorders'
$ws.Range('D2').Value = 'A negative family history reduces the probability of an inherited autoimmune process, making CREST syndrome less likely.'
$ws.Range('A3').Value = 'Family history of scleroderma or CREST syndrome'
$ws.Range('B3').Value = 'Scleroderma and its variants, including CREST syndrome, have genetic predispositions; the presence of affected family members is a supportive finding.'
$ws.Range('C3').Value = 'Family history predominantly indicates gastrointestinal disorders of non-autoimmune origin'
$ws.Range('D3').Value = 'If familial gastrointestinal issues are present but without accompanying autoimmune or connective tissue pathology, this lowers the suspicion for CREST syndrome–related esophageal dysmotility.'
$ws.Range('A4').Value = 'Family history of esophageal dysmotility or unexplained gastrointestinal motility issues'
$ws.Range('B4').Value = 'Given that Type 2 Achalasia presents with esophageal motility problems, a family history indicating similar gastrointestinal dysmotility supports a genetic predisposition to disorders affecting the esophagus.'
$ws.Range('C4').Value = 'Social environment with known exposure to environmental factors linked to more common causes of esophageal dysmotility (e.g., lifestyle-related gastroesophageal reflux)'
$ws.Range('D4').Value = 'Lifestyle factors that contribute to esophageal motility disorders other than autoimmune conditions (like reflux disease) argue against a diagnosis of CREST syndrome with Type 2 Achalasia.'
$ws.Range('A5').Value = 'Personal background indicating early-onset Raynaud phenomenon (in family context)'
$ws.Range('B5').Value = 'Raynaud phenomenon is a common and early feature of CREST syndrome; a familial tendency towards this phenomenon indirectly supports the likelihood of an underlying connective tissue disorder including CREST.'
$ws.Range('C5').Value = 'Lack of any family history of Raynaud phenomenon or scleroderma'
$ws.Range('D5').Value = 'The absence of classical features such as Raynaud phenomenon or scleroderma in the family decreases the likelihood of a familial, autoimmune variant like CREST syndrome.'
$ws.Range('A6').Value = 'Ethnic background with reported increased prevalence of autoimmune connective tissue disorders'
$ws.Range('B6').Value = 'Some ethnic groups have a higher predisposition to autoimmune conditions; belonging to such groups bolsters the relative likelihood of developing CREST syndrome.'
$ws.Range('C6').Value = 'Social history indicating robust overall health with minimal exposure to stressors linked to autoimmune triggers'
$ws.Range('D6').Value = 'A social history marked by low exposure to potential autoimmune triggers and a generally healthy family background argues against the development of a systemic autoimmune condition such as CREST syndrome.'

# Sheet: obj
$ws = $wb.Worksheets.Item('obj')
$ws.Range('A2').Value = 'Skin thickening of the fingers and face.'
$ws.Range('B2').Value = 'This is a classic sign of sclerodactyly found in CREST syndrome and supports exposure to systemic sclerosis features including esophageal dysmotility seen in Type 2 Achalasia.'
$ws.Range('C2').Value = 'Normal skin texture without evidence of thickening or tightening.'
$ws.Range('D2').Value = 'Lack of skin thickening makes the diagnosis of scleroderma-related CREST syndrome less likely, as skin involvement is a hallmark of the condition.'
$ws.Range('A3').Value = 'Presence of telangiectasias on the face and nailfolds.'
$ws.Range('B3').Value = 'Telangiectasias are a defining feature of CREST syndrome and their presence on physical exam helps distinguish the condition with high specificity.'
$ws.Range('C3').Value = 'Absence of telangiectasias on the face and nailfolds.'
$ws.Range('D3').Value = 'Not detecting telangiectasias decreases the likelihood of CREST syndrome, which typically exhibits these vascular changes.'
$ws.Range('A4').Value = 'Reduced skin elasticity and tight, shiny skin texture.'
$ws.Range('B4').Value = 'The decreased pliability of the skin is indicative of fibrotic changes in systemic sclerosis, supporting the CREST diagnosis.'
$ws.Range('C4').Value = 'Intact nailfold capillary structure without dilated or dropout capillaries.'
$ws.Range('D4').Value = 'An unremarkable nailfold capillary exam is not consistent with systemic sclerosis, reducing the probability of CREST syndrome.'
$ws.Range('A5').Value = 'Digital pitting scars and sclerodactyly.'
$ws.Range('B5').Value = 'Digital pitting is often seen in CREST syndrome and points to microvascular damage, reinforcing the presence of systemic sclerosis.'
$ws.Range('C5').Value = 'No evidence of digital pitting, ulceration, or sclerodactyly.'
$ws.Range('D5').Value = 'The absence of these digital abnormalities argues against the diagnosis of CREST syndrome, a disease known for such vascular and fibrotic changes.'
$ws.Range('A6').Value = 'Observed color changes in the digits consistent with Raynaud''s phenomenon upon temperature challenge.'
$ws.Range('B6').Value = 'Raynaud''s phenomenon is frequently associated with CREST syndrome and, when triggered during exam, offers supportive evidence for the diagnosis.'
$ws.Range('C6').Value = 'Lack of color changes in the digits during cold challenge testing.'
$ws.Range('D6').Value = 'The absence of Raynaud''s phenomenon on physical exam significantly reduces the likelihood of CREST syndrome, as this is a key diagnostic feature.'

# Sheet: test
$ws = $wb.Worksheets.Item('test')
$ws.Range('A2').Value = 'High-resolution esophageal manometry showing pan-esophageal pressurization'
$ws.Range('B2').Value = 'Pan-esophageal pressurization is a hallmark finding in Type II Achalasia and strongly supports the diagnosis when present, as it occurs much more frequently in these patients compared to controls.'
$ws.Range('C2').Value = 'Normal esophageal manometry with preserved peristalsis and normal lower esophageal sphincter relaxation'
$ws.Range('D2').Value = 'A normal manometric study effectively rules out achalasia, including Type II, thereby arguing strongly against CREST syndrome with Type II Achalasia.'
$ws.Range('A3').Value = 'Barium swallow demonstrating retained contrast with a smooth tapered narrowing at the gastroesophageal junction (bird-beak appearance)'
$ws.Range('B3').Value = 'The classic ''bird-beak'' appearance associated with achalasia, particularly with contrast retention, suggests impaired esophageal emptying, which aligns with Type II Achalasia seen in CREST syndrome.'
$ws.Range('C3').Value = 'Barium swallow demonstrating normal esophageal emptying without a ''bird-beak'' deformity'
$ws.Range('D3').Value = 'The absence of a ''bird-beak'' narrowing and normal contrast passage on a barium study makes a diagnosis of achalasia unlikely, reducing the likelihood of concurrent CREST syndrome affecting the esophagus.'
$ws.Range('A4').Value = 'Positive serum anti-centromere antibodies'
$ws.Range('B4').Value = 'Anti-centromere antibodies are highly specific for limited cutaneous systemic sclerosis (CREST syndrome) and their presence in combination with achalasia findings provides strong evidence favoring the diagnosis.'
$ws.Range('C4').Value = 'Negative serum anti-centromere antibody test'
$ws.Range('D4').Value = 'Since anti-centromere antibodies are closely associated with CREST syndrome, their absence serves as strong evidence against the diagnosis in the appropriate clinical context.'
$ws.Range('A5').Value = 'Esophageal pH monitoring showing minimal acid exposure despite dysphagia'
$ws.Range('B5').Value = 'The finding of normal acid exposure in a dysmotile esophagus supports a primary motility disorder (Type II Achalasia) rather than reflux-related changes, helping to differentiate CREST-associated esophageal dysmotility from GERD.'
$ws.Range('C5').Value = 'Normal esophageal pH monitoring indicating typical acid exposure patterns'
$ws.Range('D5').Value = 'A normal pH study without abnormal motility-related reflux changes weakens the evidence for a motility disorder like Type II Achalasia associated with CREST syndrome.'
$ws.Range('A6').Value = 'High-resolution chest CT indicating esophageal dilation with subtle interstitial changes'
$ws.Range('B6').Value = 'Chest CT findings that demonstrate esophageal dilation alongside interstitial alterations can associate with CREST syndrome-related connective tissue involvement and impaired esophageal function seen in Type II Achalasia.'
$ws.Range('C6').Value = 'High-resolution chest CT showing a normal-caliber esophagus with no evidence of dilation or interstitial changes'
$ws.Range('D6').Value = 'A CT scan that reveals a structurally normal esophagus and lacks features of connective tissue involvement would argue strongly against the presence of CREST syndrome with esophageal dysmotility.'
